$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'70.347.13"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.13%  "
$ws.Range("D3").Value = "'3.550.91"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.06%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'618.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.62%  "
$ws.Range("D6").Value = "'187.80"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.64%  "
$ws.Range("D7").Value = "'0.638"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.57%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -0.55%  "
$ws.Range("E10").Value = "  +1.77%  "
$ws.Range("D11").Value = "'53.92"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.01%  "
$ws.Range("E12").Value = "  -4.13%  "
$ws.Range("D13").Value = "'9.72"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.68%  "
$ws.Range("D14").Value = "'4.117.83"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.85%  "
$ws.Range("D15").Value = "'617.15"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +7.63%  "
$ws.Range("D16").Value = "'70.378.14"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.88%  "
$ws.Range("E17").Value = "  +3.64%  "
$ws.Range("D18").Value = "'19.16"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.87%  "
$ws.Range("D19").Value = "'3.544.30"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.46%  "
$ws.Range("E20").Value = "  +0.23%  "
$ws.Range("E21").Value = "  -1.12%  "
$ws.Range("D22").Value = "'17.69"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.32%  "
$ws.Range("D23").Value = "'105.08"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +10.46%  "
$ws.Range("E24").Value = "  +2.64%  "
$ws.Range("D25").Value = "'5.11"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("E26").Value = "  +3.52%  "
$ws.Range("D27").Value = "'11.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.14%  "
$ws.Range("D28").Value = "'9.94"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +8.80%  "
$ws.Range("D29").Value = "'34.48"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.48%  "
$ws.Range("D30").Value = "'7.11"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.75%  "
$ws.Range("D31").Value = "'12.55"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.12%  "
$ws.Range("E32").Value = "  +2.02%  "
$ws.Range("D33").Value = "'64.33"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.14%  "
$ws.Range("D34").Value = "'3.73"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +15.64%  "
$ws.Range("B35").Value = "Bittensor"
$ws.Range("C35").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D35").Value = "'538.47"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.07%  "
$ws.Range("B36").Value = "Fetch.AI"
$ws.Range("C36").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D36").Value = "'3.17"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.88%  "
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("E38").Value = "  -3.09%  "
$ws.Range("D39").Value = "'37.36"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.66%  "
$ws.Range("E40").Value = "  -3.41%  "
$ws.Range("D41").Value = "'3.551.92"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.36%  "
$ws.Range("E42").Value = "  +3.77%  "
$ws.Range("E43").Value = "  +1.97%  "
$ws.Range("E44").Value = "  +5.60%  "
$ws.Range("D45").Value = "'2.97"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.20%  "
$ws.Range("D46").Value = "'0.144"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.70%  "
$ws.Range("D47").Value = "'3.40"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.38%  "
$ws.Range("D48").Value = "'9.04"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.85%  "
$ws.Range("E49").Value = "  +0.31%  "
$ws.Range("D50").Value = "'134.54"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.26%  "
$ws.Range("D51").Value = "'1.38"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.65%  "
